$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Create new row 28 by cloning the formatting of row 27 (same visual
#    pattern the template uses for every new "group" row in this block).
# ---------------------------------------------------------------------------
$ws.Range("A27:O27").Copy()
$ws.Range("A28:O28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Stash the current formatting of the cells whose merge-area is about to
#    change, so it can be restored after Merge()/UnMerge() - Excel's merge
#    operation renormalizes the borders of the merged block and would
#    otherwise silently create new style records.
# ---------------------------------------------------------------------------
$ws.Range("A25").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("C25").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("D25").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("E27").Copy()
$ws.Range("Z5").PasteSpecial(-4122)
$ws.Range("F27").Copy()
$ws.Range("Z6").PasteSpecial(-4122)
$ws.Range("G27").Copy()
$ws.Range("Z7").PasteSpecial(-4122)
$ws.Range("H27").Copy()
$ws.Range("Z8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Re-shape the merged blocks to include the new row 28.
# ---------------------------------------------------------------------------
$ws.Range("A25:A27").UnMerge()
$ws.Range("A25:A28").Merge()

$ws.Range("B25:B27").UnMerge()
$ws.Range("B25:B28").Merge()

$ws.Range("C25:C27").UnMerge()
$ws.Range("C25:C28").Merge()

$ws.Range("D25:D27").UnMerge()
$ws.Range("D25:D28").Merge()

$ws.Range("E27:E28").Merge()
$ws.Range("F27:F28").Merge()
$ws.Range("G27:G28").Merge()
$ws.Range("H27:H28").Merge()

# ---------------------------------------------------------------------------
# 4. Restore the original formatting that Merge()/UnMerge() perturbed.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("A25:A28").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B25:B28").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("C25:C28").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("D25:D28").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("E27:E28").PasteSpecial(-4122)
$ws.Range("Z6").Copy()
$ws.Range("F27:F28").PasteSpecial(-4122)
$ws.Range("Z7").Copy()
$ws.Range("G27:G28").PasteSpecial(-4122)
$ws.Range("Z8").Copy()
$ws.Range("H27:H28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("Z1:Z8").Clear()

# ---------------------------------------------------------------------------
# 5. Change I27 / J27 formatting (style only, value/text unchanged) to match
#    the "I8 / J8" look used elsewhere for this column pair.
# ---------------------------------------------------------------------------
$ws.Range("I8").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("J8").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Apply the numeric value changes from the R script re-run.
# ---------------------------------------------------------------------------
$ws.Range("K25").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = 154996
$ws.Range("K27").Value = 3354

# ---------------------------------------------------------------------------
# 7. Populate the new row 28 with its data. I28 holds a numeric-looking code
#    that must stay text (as in the rest of the "Natureza Despesa" column),
#    so build it with a TEXT() formula and immediately flatten it back to a
#    plain cached value/shared-string, keeping the cell's existing style.
# ---------------------------------------------------------------------------
$ws.Range("I28").Formula = '=TEXT(339018,"0")'
$ws.Range("I28").Copy()
$ws.Range("I28").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("J28").Value = "AUXILIO FINANCEIRO A ESTUDANTES"
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = 5641
